$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. First paragraph: "This is a Microsoft word document." gets a
#    trailing two spaces (kept in its own un-colored run) followed by
#    three new red-colored runs spelling out
#    "(This is a change – Version for main branch)"
#    split exactly as in the target XML:
#      run A: "(This is a change – Ve"
#      run B: "rsion for main branch"
#      run C: ")"
# ------------------------------------------------------------------

$d.Content.Find.Execute(
    "This is a Microsoft word document.", $true, $false, $false, $false,
    $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

$p1 = $d.Paragraphs(1).Range
$insertPos = $p1.End - 1

$segments = @(
    "(This is a change " + [char]0x2013 + " Ve",
    "rsion for main branch",
    ")"
)

foreach ($seg in $segments) {
    $ip = $d.Range($insertPos, $insertPos)
    $ip.InsertAfter($seg)
    $segEnd = $insertPos + $seg.Length
    $segRange = $d.Range($insertPos, $segEnd)
    $segRange.Font.Color = 255
    $insertPos = $segEnd
}

# ------------------------------------------------------------------
# 2. Remove the final paragraph in the body ("ank God almighty, we
#    are free at last.") entirely, including its paragraph mark, so
#    the previous paragraph ("            Shall be lifted-nevermore!")
#    becomes the document's last paragraph.
# ------------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastPara.Range.Delete()

# ------------------------------------------------------------------
# 3. Prune the handful of unused custom styles that the resaved
#    document no longer carries (Heading 2/4, the scraped-blog
#    character/paragraph styles, etc.). None of them are referenced
#    by any paragraph/run left in the document.
#
#    NOTE: Styles.Item(<name>) resolves against each style's
#    *original* position, so deleting in ascending order shifts the
#    backing store underneath later lookups. Deleting from the
#    highest original index down avoids that.
# ------------------------------------------------------------------

$styleNamesHighToLow = @(
    "podcast-tools__subscribe-links",
    "generic-title",
    "subscribe-more-info",
    "subscribe",
    "audio-tool",
    "Heading 4 Char",
    "Heading 2 Char",
    "Hyperlink",
    "apple-converted-space",
    "Heading 4",
    "Heading 2"
)

foreach ($styleName in $styleNamesHighToLow) {
    $style = $d.Styles.Item($styleName)
    $style.Delete()
}
